$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Drop the trailing three slides (old slide6 "User Engagement Patterns",
#    slide7 "Application Usage Analysis", slide8 "Key Recommendations").
#    Their content has effectively been folded into the new slide 5
#    ("Recommendations"), so the deck goes from 8 slides down to 5.
# ---------------------------------------------------------------------------
$p.Slides.Item(8).Delete()
$p.Slides.Item(7).Delete()
$p.Slides.Item(6).Delete()

# ---------------------------------------------------------------------------
# 2) Slide 1 - title slide: refresh the subtitle copy.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Interim Findings"

# ---------------------------------------------------------------------------
# Slides 2 and 5 both started life with a plain "Text Placeholder" body
# (inherited from the Section Header layout) that never got any content.
# They need a real Content Placeholder instead, so clone the fully-formed
# one from slide 3 (still pristine at this point) onto both slides first,
# before any of slide 3/4's own bullet edits happen below.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s3 = $p.Slides.Item(3)
$s4 = $p.Slides.Item(4)
$s5 = $p.Slides.Item(5)

$p.Slides.Item(3).Shapes.Item(2).Copy()
$newBody2 = $s2.Shapes.Paste().Item(1)
$s2.Shapes.Item(2).Delete()

$p.Slides.Item(3).Shapes.Item(2).Copy()
$newBody5 = $s5.Shapes.Paste().Item(1)
$s5.Shapes.Item(2).Delete()

# ---------------------------------------------------------------------------
# 3) Slide 2 - was "Task 1: User Overview Analysis"; becomes "Key Findings
#    and Growth Opportunities" with five bullet points.
# ---------------------------------------------------------------------------
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Key Findings and Growth Opportunities"
$newBody2.TextFrame.TextRange.Text = "`rHigh-end smartphones dominate usage`rSocial media and video streaming drive engagement`rGaming apps show growing usage`rPremium users are key demographic for growth`rNetwork optimization opportunities identified"

# ---------------------------------------------------------------------------
# 4) Slide 3 - "Top 10 Handsets" -> "Device Preferences"; trim the bullet
#    list from five items to four and refresh the remaining copy.
# ---------------------------------------------------------------------------
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Device Preferences"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$line = $tr3.Lines(5)
$tr3.Characters($line.Start - 1, $line.Length + 1).Delete()

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Lines(2).Text = "Top handsets dominated by premium devices"
$tr3.Lines(3).Text = "Apple leads manufacturer market share"
$tr3.Lines(4).Text = "High-end devices show increased data usage"
$tr3.Lines(5).Text = "Opportunity for targeted premium services"

# ---------------------------------------------------------------------------
# 5) Slide 4 - "Top Manufacturers" -> "User Engagement Insights"; trim the
#    bullet list from five items to four and refresh the remaining copy.
# ---------------------------------------------------------------------------
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "User Engagement Insights"

$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$line = $tr4.Lines(5)
$tr4.Characters($line.Start - 1, $line.Length + 1).Delete()

$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Lines(2).Text = "Peak usage patterns identified in evening hours"
$tr4.Lines(3).Text = "Social media drives majority of traffic"
$tr4.Lines(4).Text = "Video streaming shows high engagement"
$tr4.Lines(5).Text = "Gaming emerges as growth segment"

# ---------------------------------------------------------------------------
# 6) Slide 5 - was "Task 2: User Engagement Analysis"; becomes
#    "Recommendations" with five bullet points (folding in the old closing
#    slide's ideas).
# ---------------------------------------------------------------------------
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Recommendations"
$newBody5.TextFrame.TextRange.Text = "`rFocus on premium smartphone users`rOptimize network for video streaming`rPartner with top manufacturers`rDevelop targeted marketing campaigns`rImplement user segmentation strategies"
